$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize A23:A24 date formatting to match the rest of the date column ---
# (copy the format that the other date cells, e.g. A2, already use)
$ws.Range("A2").Copy()
$ws.Range("A23:A24").PasteSpecial(-4122)

# --- Add new row 25: "Unique Paths II" ---

# A25: date value, formatted like the other date cells in column A
$ws.Range("A2").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value = 46056

# B25: plain text question name
$ws.Range("B25").Value = "Unique Paths II"

# C25: url text, formatted + hyperlinked like the other links in column C
$ws.Range("C24").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = "https://leetcode.com/problems/unique-paths-ii/"
$ws.Hyperlinks.Add($ws.Range("C25"), "https://leetcode.com/problems/unique-paths-ii/")

# Adding the hyperlink registers a built-in "Hyperlink" cell style that this
# workbook didn't previously use; drop it again since the other hyperlinked
# URL cells in this sheet are plain (unstyled) text.
$wb.Styles.Item("Hyperlink").Delete()

# Adding the hyperlink can also reset C25's cell formatting; reapply it so the
# cell keeps matching the rest of the hyperlinked URL column.
$ws.Range("C24").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = "https://leetcode.com/problems/unique-paths-ii/"
